$d = $word.ActiveDocument

$pairs = @(
    @{old = "410÷6=68, 2"; new = "749÷9=83, 2"},
    @{old = "100÷5=20, 0"; new = "157÷2=78, 1"},
    @{old = "419÷7=59, 6"; new = "726÷6=121, 0"},
    @{old = "408÷7=58, 2"; new = "834÷8=104, 2"},
    @{old = "855÷5=171, 0"; new = "298÷6=49, 4"},
    @{old = "686÷5=137, 1"; new = "560÷6=93, 2"},
    @{old = "623÷5=124, 3"; new = "750÷4=187, 2"},
    @{old = "491÷6=81, 5"; new = "965÷2=482, 1"},
    @{old = "811÷2=405, 1"; new = "451÷8=56, 3"},
    @{old = "619÷7=88, 3"; new = "774÷6=129, 0"},
    @{old = "747÷8=93, 3"; new = "214÷8=26, 6"},
    @{old = "965÷6=160, 5"; new = "515÷7=73, 4"},
    @{old = "312÷6=52, 0"; new = "712÷5=142, 2"},
    @{old = "381÷7=54, 3"; new = "662÷4=165, 2"},
    @{old = "548÷6=91, 2"; new = "894÷8=111, 6"},
    @{old = "175÷8=21, 7"; new = "250÷2=125, 0"},
    @{old = "233÷3=77, 2"; new = "392÷8=49, 0"},
    @{old = "416÷4=104, 0"; new = "687÷8=85, 7"},
    @{old = "644÷6=107, 2"; new = "688÷7=98, 2"},
    @{old = "709÷3=236, 1"; new = "965÷2=482, 1"},
    @{old = "916÷6=152, 4"; new = "939÷5=187, 4"},
    @{old = "220÷3=73, 1"; new = "491÷9=54, 5"},
    @{old = "250÷5=50, 0"; new = "340÷4=85, 0"},
    @{old = "149÷4=37, 1"; new = "615÷4=153, 3"},
    @{old = "546÷4=136, 2"; new = "292÷3=97, 1"}
)

foreach ($pair in $pairs) {
    $range = $d.Content
    $range.Find.Execute($pair.old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $pair.new, 2)
}
